# Auto-generated edit script: update Sheets via scheduled runner
# Applies cached market-data value updates to the Leve profit tables.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 137
$ws.Range("H137").Value = 898.9722
$ws.Range("I137").Value = 798.80646
$ws.Range("K137").Value = 2396.41938
$ws.Range("M137").Value = 153.5806199999997

# ALC row 138
$ws.Range("H138").Value = 3682.606
$ws.Range("I138").Value = 1465.6154
$ws.Range("J138").Value = 4472.219
$ws.Range("K138").Value = 4396.8462
$ws.Range("L138").Value = 13416.657
$ws.Range("M138").Value = 743.1538
$ws.Range("N138").Value = -23696.657

$ws = $wb.Worksheets.Item("ARM")
# ARM row 32
$ws.Range("H32").Value = 20092.98
$ws.Range("I32").Value = 14283.694
$ws.Range("J32").Value = 91256.75
$ws.Range("K32").Value = 14283.694
$ws.Range("L32").Value = 91256.75
$ws.Range("M32").Value = -13996.694
$ws.Range("N32").Value = -91830.75

# ARM row 36
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("N36").ClearContents()

# ARM row 37
$ws.Range("H37").Value = 8202.799999999999
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 8202.799999999999
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 8202.799999999999
$ws.Range("M37").ClearContents()
$ws.Range("N37").Value = -8748.799999999999

# ARM row 110
$ws.Range("H110").Value = 630.7692
$ws.Range("I110").Value = 534.63635
$ws.Range("J110").Value = 1159.5
$ws.Range("K110").Value = 534.63635
$ws.Range("L110").Value = 1159.5
$ws.Range("M110").Value = 1510.36365
$ws.Range("N110").Value = -5249.5

$ws = $wb.Worksheets.Item("BSM")
# BSM row 86
$ws.Range("H86").Value = 23046
$ws.Range("I86").Value = 52499.5
$ws.Range("K86").Value = 52499.5
$ws.Range("M86").Value = -51376.5

# BSM row 89
$ws.Range("H89").Value = 23046
$ws.Range("I89").Value = 52499.5
$ws.Range("K89").Value = 262497.5
$ws.Range("M89").Value = -256881.5

# BSM row 105
$ws.Range("H105").Value = 4372.294
$ws.Range("I105").Value = 6002.25
$ws.Range("K105").Value = 6002.25
$ws.Range("M105").Value = -4255.25

$ws = $wb.Worksheets.Item("CRP")
# CRP row 16
$ws.Range("H16").Value = 1308.6364
$ws.Range("I16").Value = 1050
$ws.Range("J16").Value = 1998.3334
$ws.Range("K16").Value = 1050
$ws.Range("L16").Value = 1998.3334
$ws.Range("M16").Value = -763
$ws.Range("N16").Value = -2572.3334

# CRP row 35
$ws.Range("H35").Value = 769.8
$ws.Range("I35").Value = 769.8
$ws.Range("K35").Value = 769.8
$ws.Range("M35").Value = -475.8

# CRP row 62
$ws.Range("H62").Value = 6022.5386
$ws.Range("I62").Value = 6074.375
$ws.Range("K62").Value = 6074.375
$ws.Range("M62").Value = -5450.375

# CRP row 65
$ws.Range("H65").Value = 6022.5386
$ws.Range("I65").Value = 6074.375
$ws.Range("K65").Value = 30371.875
$ws.Range("M65").Value = -27251.875

# CRP row 105
$ws.Range("H105").Value = 9000
$ws.Range("I105").Value = 9000
$ws.Range("K105").Value = 9000
$ws.Range("M105").Value = -7253

# CRP row 113
$ws.Range("H113").Value = 1308.6364
$ws.Range("I113").Value = 1050
$ws.Range("J113").Value = 1998.3334
$ws.Range("K113").Value = 1050
$ws.Range("L113").Value = 1998.3334
$ws.Range("M113").Value = 1120
$ws.Range("N113").Value = -6338.3334

# CRP row 132
$ws.Range("H132").Value = 1531.4
$ws.Range("I132").Value = 1153.8823
$ws.Range("K132").Value = 3461.6469
$ws.Range("M132").Value = -931.6468999999997

$ws = $wb.Worksheets.Item("CUL")
# CUL row 99
$ws.Range("H99").Value = 254624.75
$ws.Range("I99").Value = 336166.34
$ws.Range("J99").Value = 10000
$ws.Range("K99").Value = 1008499.02
$ws.Range("L99").Value = 30000
$ws.Range("M99").Value = -1006253.02
$ws.Range("N99").Value = -34492

# CUL row 100
$ws.Range("H100").Value = 6600
$ws.Range("J100").Value = 6600
$ws.Range("L100").Value = 19800
$ws.Range("N100").Value = -21422

# CUL row 102
$ws.Range("H102").Value = 4183.3335
$ws.Range("J102").Value = 4183.3335
$ws.Range("L102").Value = 12550.0005
$ws.Range("N102").Value = -17418.0005

# CUL row 107
$ws.Range("H107").Value = 557047.1
$ws.Range("J107").Value = 973288.4
$ws.Range("L107").Value = 2919865.2
$ws.Range("N107").Value = -2923705.2

# CUL row 108
$ws.Range("H108").Value = 607.5
$ws.Range("I108").Value = 406
$ws.Range("J108").Value = 1615
$ws.Range("K108").Value = 1218
$ws.Range("L108").Value = 4845
$ws.Range("M108").Value = 1662
$ws.Range("N108").Value = -10605

# CUL row 117
$ws.Range("H117").Value = 2094.5
$ws.Range("J117").Value = 3410.25
$ws.Range("L117").Value = 10230.75
$ws.Range("N117").Value = -17114.75

# CUL row 122
$ws.Range("H122").Value = 21591
$ws.Range("I122").Value = 582.625
$ws.Range("J122").Value = 25690.195
$ws.Range("K122").Value = 5243.625
$ws.Range("L122").Value = 231211.755
$ws.Range("M122").Value = -2793.625
$ws.Range("N122").Value = -236111.755

# CUL row 129
$ws.Range("H129").Value = 27674.125
$ws.Range("I129").Value = 1252.7273
$ws.Range("J129").Value = 37696.035
$ws.Range("K129").Value = 3758.1819
$ws.Range("L129").Value = 113088.105
$ws.Range("M129").Value = 1241.8181
$ws.Range("N129").Value = -123088.105

# CUL row 134
$ws.Range("H134").Value = 4194.9614
$ws.Range("I134").Value = 2977.3157
$ws.Range("K134").Value = 8931.947100000001
$ws.Range("M134").Value = -3861.947100000001

# CUL row 139
$ws.Range("H139").Value = 35993.484
$ws.Range("I139").Value = 42826.293
$ws.Range("J139").Value = 3196
$ws.Range("K139").Value = 128478.879
$ws.Range("L139").Value = 9588
$ws.Range("M139").Value = -123338.879
$ws.Range("N139").Value = -19868

# CUL row 140
$ws.Range("H140").Value = 276982.72
$ws.Range("I140").Value = 433630
$ws.Range("J140").Value = 2850
$ws.Range("K140").Value = 1300890
$ws.Range("L140").Value = 8550
$ws.Range("M140").Value = -1295710
$ws.Range("N140").Value = -18910

$ws = $wb.Worksheets.Item("GSM")
# GSM row 132
$ws.Range("H132").Value = 2741.5789
$ws.Range("I132").Value = 2349.2666
$ws.Range("J132").Value = 4212.75
$ws.Range("K132").Value = 7047.7998
$ws.Range("L132").Value = 12638.25
$ws.Range("M132").Value = -4517.7998
$ws.Range("N132").Value = -17698.25

$ws = $wb.Worksheets.Item("LTW")
# LTW row 7
$ws.Range("H7").Value = 2593.625
$ws.Range("I7").Value = 2432.5334
$ws.Range("K7").Value = 2432.5334
$ws.Range("M7").Value = -2320.5334

# LTW row 32
$ws.Range("H32").Value = 1006.5
$ws.Range("I32").Value = 1006.5
$ws.Range("K32").Value = 1006.5
$ws.Range("M32").Value = -689.5

# LTW row 46
$ws.Range("H46").Value = 25000688
$ws.Range("I46").Value = 28572142
$ws.Range("J46").Value = 500
$ws.Range("K46").Value = 28572142
$ws.Range("L46").Value = 500
$ws.Range("M46").Value = -28571954
$ws.Range("N46").Value = -876

# LTW row 50
$ws.Range("H50").Value = 8000
$ws.Range("J50").Value = 8000
$ws.Range("L50").Value = 8000
$ws.Range("N50").Value = -9274

# LTW row 126
$ws.Range("H126").Value = 2593.625
$ws.Range("I126").Value = 2432.5334
$ws.Range("K126").Value = 7297.600199999999
$ws.Range("M126").Value = -4827.600199999999

$ws = $wb.Worksheets.Item("WVR")
# WVR row 62
$ws.Range("H62").Value = 7971.4287
$ws.Range("I62").Value = 6600
$ws.Range("J62").Value = 9000
$ws.Range("K62").Value = 6600
$ws.Range("L62").Value = 9000
$ws.Range("M62").Value = -5976
$ws.Range("N62").Value = -10248

# WVR row 65
$ws.Range("H65").Value = 7971.4287
$ws.Range("I65").Value = 6600
$ws.Range("J65").Value = 9000
$ws.Range("K65").Value = 33000
$ws.Range("L65").Value = 45000
$ws.Range("M65").Value = -29880
$ws.Range("N65").Value = -51240

# WVR row 132
$ws.Range("H132").Value = 2303
$ws.Range("I132").Value = 2451
$ws.Range("J132").Value = 2195.3635
$ws.Range("K132").Value = 7353
$ws.Range("L132").Value = 6586.0905
$ws.Range("M132").Value = -4823
$ws.Range("N132").Value = -11646.0905

Write-Output "Applied Bahamut_Profits market data refresh"